# Updates crypto price/volume figures (and re-ranks two coin pairs) per the
# Tue Aug 29 2023 GitHub Actions refresh of the cryptos list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds plain text (values like "26.077.62" use "." as a
# thousands separator and are not real numbers). Cells whose new text DOES
# parse as a number (e.g. "1.003") are written with a leading apostrophe so
# Excel stores them as literal text instead of silently converting them to a
# Double and dropping significant trailing zeros - exactly what a manual
# quote-prefixed entry in Excel's UI would do.

$ws.Range("D2").Value = "26.077.62"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.647.52"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'218.32"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "'0.5189"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'0.2624"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").Value = "'0.06298"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "'20.25"
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("D11").Value = "'0.07688"
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").Value = "'4.589"
$ws.Range("E12").Value = "  +2.53%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.875.55"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.607.15"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").Value = "'0.5565"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "0.0₅8107"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").Value = "'65.06"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "26.062.73"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "'4.604"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'192.77"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "'10.43"
$ws.Range("E22").Value = "  +3.29%  "
$ws.Range("D23").Value = "'5.912"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").Value = "'143.98"
$ws.Range("E25").Value = "  -1.68%  "
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("D27").Value = "'7.181"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").Value = "'1.510"
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("D30").Value = "'0.05340"
$ws.Range("E30").Value = "  -4.96%  "
$ws.Range("D31").Value = "'1.268"
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").Value = "'3.323"
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("D35").Value = "'2.419"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").Value = "'2.781"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").Value = "'0.9398"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").Value = "'0.5591"
$ws.Range("D39").Value = "'0.01570"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").Value = "'5.773"
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "1.025.20"
$ws.Range("E42").Value = "  -2.49%  "
$ws.Range("D43").Value = "'0.8240"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").Value = "'100.86"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").Value = "1.786.00"
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  +6.10%  "
$ws.Range("D47").Value = "'57.29"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").Value = "'0.4311"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").Value = "'7.898"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").Value = "'0.05109"
$ws.Range("E51").Value = "  -3.87%  "
